$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "Nombre"
$ws.Range("C1").Value = "Teléfono"
$ws.Range("D1").Value = "EMAIL"

# --- Data rows: Name, Phone, concatenated id formula, and an "X" marker column ---
$names  = @("Rodrigo", "Rubén", "Lilia", "Diego")
$phones = @(5545125300, 5522654896, 5510721906, 5526547119)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $names[$i]
    $ws.Cells.Item($row, 3).Value = $phones[$i]
    $ws.Cells.Item($row, 1).Formula = "=CONCAT(B" + $row + ",C" + $row + ")"
    $ws.Cells.Item($row, 4).Value = "X"
}

# Apply a style to the Name/Phone cells
$ws.Range("B2:C5").Style = "Normal"

# Center the "X" marker column with its own style
$ws.Range("D2:D5").HorizontalAlignment = -4108

# Widen column A ("id") so the concatenated name+phone text is visible
$ws.Columns.Item(1).ColumnWidth = 27.826666666666668

$ws.Range("F9").Select() | Out-Null
